{"js": "// Translate the remaining English strings in HIVE TEAMS.docx into German.\nconst replacements = [\n  {\n    find: \"This Hive is responsible for on-boarding & generalized SmartCash support.\",\n    replace: \"Dieser Hive ist verantwortlich f\u00fcr On-Boarding & allgemeinen SmartCash Support.\"\n  },\n  {\n    find: \"Hive Coordinator\",\n    replace: \"Hive Koordinator\"\n  },\n  {\n    find: \"Alex is a jack of all trades who loves Technology, Graphics, Web Design & Infrastructure.\",\n    replace: \"Alex ist ein Allesk\u00f6nner, der Technologie, Grafik, Web Design & Infrastruktur liebt.\"\n  },\n  {\n    find: \"Fiscal Officer\",\n    replace: \"Steuerbeauftragter\"\n  },\n  {\n    find: \"Vice Coordinator\",\n    replace: \"Vize-Koordinator\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Translate the remaining English strings in HIVE TEAMS.docx into German.\n$d = $word.ActiveDocument\n\n# Straight text swaps: walk the paragraphs and overwrite the Range.Text of\n# any paragraph whose text matches one of the English source strings.\n$directReplacements = @(\n    @{ Find = \"This Hive is responsible for on-boarding & generalized SmartCash support.\"; Replace = \"Dieser Hive ist verantwortlich f\u00fcr On-Boarding & allgemeinen SmartCash Support.\" },\n    @{ Find = \"Hive Coordinator\"; Replace = \"Hive Koordinator\" },\n    @{ Find = \"Alex is a jack of all trades who loves Technology, Graphics, Web Design & Infrastructure.\"; Replace = \"Alex ist ein Allesk\u00f6nner, der Technologie, Grafik, Web Design & Infrastruktur liebt.\" }\n)\n\nforeach ($p in $d.Paragraphs) {\n    $r = $p.Range\n    $text = $r.Text\n    if ($text -eq $null) { continue }\n    $trimmed = $text.TrimEnd(\"`r\", \"`a\")\n    foreach ($rep in $directReplacements) {\n        if ($trimmed -eq $rep.Find) {\n            $r.Text = $rep.Replace\n            break\n        }\n    }\n}\n\n# Find & Replace for the remaining two short labels.\n$findReplacements = @(\n    @{ Find = \"Fiscal Officer\"; Replace = \"Steuerbeauftragter\" },\n    @{ Find = \"Vice Coordinator\"; Replace = \"Vize-Koordinator\" }\n)\n\nforeach ($rep in $findReplacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    [void]$find.Execute(\n        $rep.Find,   # FindText\n        $true,       # MatchCase\n        $true,       # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $rep.Replace,# ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
